$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "liefern",
    "schneiden",
    "wirken",
    "machen",
    "fangen",
    "doppeln",
    "herrschen",
    "landen",
    "schmecken",
    "tollen",
    "steuern",
    "filmen",
    "stehlen",
    "enden",
    "werden",
    "lügen",
    "gelten",
    "nutzen",
    "decken",
    "rasen",
    "jubeln",
    "scheinen",
    "drehen",
    "nennen",
    "schenken",
    "fließen",
    "formen",
    "töten",
    "jagen",
    "führen",
    "schätzen",
    "knarren",
    "brauchen",
    "kosten",
    "schauen",
    "schleppen",
    "starren",
    "bluten",
    "schwören",
    "mögen",
    "stechen",
    "arten",
    "wellen",
    "wüten",
    "platzen",
    "mühen",
    "gründen",
    "pflegen",
    "wundern",
    "klettern",
    "opfern",
    "fehlen",
    "sorgen",
    "wachsen",
    "rühren",
    "klingen",
    "seufzen",
    "retten",
    "zögern",
    "siegen",
    "leugnen",
    "bitten",
    "betteln",
    "wehren",
    "ehren",
    "pflanzen",
    "scheitern",
    "werfen",
    "zeigen",
    "ächzen",
    "dienen",
    "zahlen",
    "boxen",
    "pfeifen",
    "schlucken",
    "backen",
    "haben",
    "morden",
    "tropfen",
    "hassen",
    "schaden",
    "suchen",
    "handeln",
    "weichen",
    "garen",
    "spielen",
    "messen",
    "streifen",
    "stillen",
    "schwächen",
    "feiern",
    "loben",
    "runden",
    "spinnen",
    "trotzen",
    "foltern",
    "greifen",
    "sitzen",
    "hören",
    "lesen",
    "flüstern",
    "fallen",
    "segnen",
    "rauchen",
    "merken",
    "sichern",
    "hoffen",
    "fahren",
    "räumen",
    "liegen",
    "proben",
    "spüren",
    "quälen",
    "warnen",
    "sinken",
    "eignen",
    "kichern",
    "stammen",
    "helfen",
    "atmen",
    "schrecken",
    "trennen",
    "dringen",
    "kümmern",
    "fragen",
    "mauern",
    "stecken",
    "geben",
    "schreiten",
    "dauern",
    "lieben",
    "fällen",
    "hauen",
    "laufen",
    "bergen",
    "plaudern",
    "beten",
    "wehtun",
    "fügen",
    "irren",
    "heilen",
    "folgen",
    "fischen",
    "trauen",
    "reisen",
    "beißen",
    "zielen",
    "erben",
    "weigern",
    "sterben",
    "heulen",
    "streichen",
    "ändern",
    "gnaden",
    "planen",
    "lächeln",
    "bauen",
    "schulden",
    "schwingen",
    "bellen",
    "schließen",
    "graben",
    "zünden",
    "münzen",
    "stellen",
    "heben",
    "biegen",
    "rufen",
    "sprengen",
    "altern",
    "ärgern",
    "husten",
    "treiben",
    "wenden",
    "malen",
    "kehren",
    "äußern",
    "grüßen",
    "öffnen",
    "achten",
    "orten",
    "sperren",
    "lohnen",
    "saufen",
    "flehen",
    "lehnen",
    "lockern",
    "reizen",
    "kaufen",
    "freuen",
    "warten",
    "flüchten"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

